# Update "Pais" (countries) sheet with refreshed COVID-19 stats and re-sort
# a few rows whose "Casos totales" (column B) crossed each other, which
# changes which country label sits on which row (Excel keeps data sorted
# descending by total cases).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: timestamp footer update
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 12:31"

# Row 7: India (data refresh)
$ws.Range("B7").Value = 457621
$ws.Range("C7").Value = 1506
$ws.Range("D7").Value = 259092
$ws.Range("E7").Value = 184029
$ws.Range("G7").Value = 17
$ws.Range("H7").Value = 14500

# Row 26: Suecia -> Belgica (rows swap order / data refresh)
$ws.Range("A26").Value = "Belgica"
$ws.Range("B26").Value = 60898
$ws.Range("C26").Value = 88
$ws.Range("D26").Value = 16771
$ws.Range("E26").Value = 34405
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 9722

# Row 27: Belgica -> Suecia (rows swap order / data refresh)
$ws.Range("A27").Value = "Suecia"
$ws.Range("B27").Value = 60837
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("H27").Value = 5161

# Row 32: Indonesia (data refresh)
$ws.Range("B32").Value = 49009
$ws.Range("C32").Value = 1113
$ws.Range("D32").Value = 19658
$ws.Range("E32").Value = 26778
$ws.Range("G32").Value = 38
$ws.Range("H32").Value = 2573

# Row 40: Polonia -> Oman (rows swap order / data refresh)
$ws.Range("A40").Value = "Oman"
$ws.Range("B40").Value = 33536
$ws.Range("C40").Value = 1142
$ws.Range("D40").Value = 17972
$ws.Range("E40").Value = 15422
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 142

# Row 41: Oman -> Polonia (rows swap order / data refresh)
$ws.Range("A41").Value = "Polonia"
$ws.Range("B41").Value = 32821
$ws.Range("C41").Value = 294
$ws.Range("D41").Value = 18134
$ws.Range("E41").Value = 13291
$ws.Range("G41").Value = 21
$ws.Range("H41").Value = 1396

# Row 49: Rumania (data refresh)
$ws.Range("B49").Value = 24826
$ws.Range("C49").Value = 321
$ws.Range("D49").Value = 17391
$ws.Range("E49").Value = 5880
$ws.Range("G49").Value = 16
$ws.Range("H49").Value = 1555

# Row 54: Kazajistan (data refresh)
$ws.Range("D54").Value = 11585
$ws.Range("E54").Value = 7046

# Row 67: Chequia -> Marruecos (rows swap order / data refresh)
$ws.Range("A67").Value = "Marruecos"
$ws.Range("B67").Value = 10693
$ws.Range("C67").Value = 349
$ws.Range("D67").Value = 8426
$ws.Range("E67").Value = 2053
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 214

# Row 68: Marruecos -> Chequia (rows swap order / data refresh)
$ws.Range("A68").Value = "Chequia"
$ws.Range("B68").Value = 10651
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = 7559
$ws.Range("E68").Value = 2752
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 340

# Row 72: Malasia (data refresh)
$ws.Range("B72").Value = 8596
$ws.Range("C72").Value = 6
$ws.Range("D72").Value = 8231
$ws.Range("E72").Value = 244

# Row 75: Finlandia (data refresh)
$ws.Range("B75").Value = 7167
$ws.Range("C75").Value = 12
$ws.Range("E75").Value = 440

# Row 77: Senegal (data refresh)
$ws.Range("B77").Value = 6129
$ws.Range("C77").Value = 95
$ws.Range("D77").Value = 4072
$ws.Range("E77").Value = 1964
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 93

# Row 106: Albania (data refresh)
$ws.Range("B106").Value = 2114
$ws.Range("C106").Value = 67
$ws.Range("D106").Value = 1217
$ws.Range("E106").Value = 850
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 47

# Row 124: Hong Kong (data refresh)
$ws.Range("B124").Value = 1180
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 1086
$ws.Range("E124").Value = 88

# Row 202: Dominica -> Fiyi (rows swap order / data refresh)
$ws.Range("A202").Value = "Fiyi"

# Row 203: Fiyi -> Dominica (rows swap order / data refresh)
$ws.Range("A203").Value = "Dominica"

# Row 211: Montserrat -> Seychelles (rows swap order / data refresh)
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 212: Seychelles -> Montserrat (rows swap order / data refresh)
$ws.Range("A212").Value = "Montserrat"
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
